$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2911
$ws.Range("K3").Value = 2849
$ws.Range("K4").Value = 590
$ws.Range("K5").Value = 183
$ws.Range("K6").Value = 3458
$ws.Range("K7").Value = 9991

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 189
$ws.Range("K4").Value = 36
$ws.Range("K6").Value = 218
$ws.Range("K7").Value = 653

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 218

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 110
$ws.Range("K7").Value = 396

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 87
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 330

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 298
$ws.Range("K8").Value = 653
$ws.Range("K11").Value = 210
$ws.Range("K12").Value = 16
$ws.Range("K19").Value = 299
$ws.Range("K20").Value = 231
$ws.Range("K23").Value = 86
$ws.Range("K27").Value = 103
$ws.Range("K29").Value = 523
$ws.Range("K33").Value = 396
$ws.Range("K34").Value = 48
$ws.Range("K37").Value = 330
$ws.Range("K42").Value = 348
$ws.Range("K43").Value = 89
$ws.Range("K44").Value = 95
$ws.Range("K46").Value = 21
$ws.Range("K51").Value = 112
$ws.Range("K52").Value = 282
$ws.Range("K54").Value = 188
$ws.Range("K65").Value = 234
$ws.Range("K66").Value = 34
$ws.Range("K70").Value = 17
$ws.Range("K72").Value = 46
$ws.Range("K76").Value = 154
$ws.Range("K79").Value = 253
$ws.Range("K83").Value = 218
$ws.Range("K85").Value = 476
$ws.Range("K94").Value = 119
$ws.Range("K96").Value = 135
$ws.Range("K98").Value = 57
$ws.Range("K101").Value = 9991

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 127

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 59
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 188

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 144
$ws.Range("K3").Value = 174
$ws.Range("K4").Value = 30
$ws.Range("K7").Value = 523

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 96
$ws.Range("K3").Value = 81
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 299

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 27
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 89
$ws.Range("K3").Value = 111
$ws.Range("K6").Value = 133
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 8
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 253

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 81
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 231

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 101
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 32
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K2").Value = 11
$ws.Range("K6").Value = 34

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 64
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 210

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 17

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 172
$ws.Range("K3").Value = 163
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 476

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K3").Value = 14
$ws.Range("K6").Value = 46

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 76
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 16
